# Insert a new weekly price record as row 529 in the Espinaca sheet.
# This shifts the existing rows 529:578 down to 530:579 (preserving all
# of their data/formatting), and populates the newly created row 529
# with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 529, pushing
# everything from 529 downward (formatting included).
$ws.Rows.Item(529).Insert()

# Populate the new row 529 with the new data point.
$ws.Range("A529").Value = 3
$ws.Range("B529").Value = "Femacal de La Calera"
$ws.Range("C529").Value = "Coquimbo"
$ws.Range("D529").Value = 45132
$ws.Range("E529").Value = 5
$ws.Range("F529").Value = 100112012
$ws.Range("G529").Value = "Espinaca"
$ws.Range("H529").Value = "Sin especificar"
$ws.Range("I529").Value = "Primera"
$ws.Range("J529").Value = 65
$ws.Range("K529").Value = 4500
$ws.Range("L529").Value = 4500
$ws.Range("M529").Value = 4500
$ws.Range("N529").Value = "$/docena de atados (3 kilos)"
$ws.Range("O529").Value = "Provincia de Quillota"
$ws.Range("P529").Value = 1500
$ws.Range("Q529").Value = 3
$ws.Range("R529").Value = "Hortaliza"
